# Populate Sheet1 with the "twoStep" results table (15 rows x 3 cols).
# Columns B, C, and A mostly hold numeric-looking values that must be
# stored as TEXT (quote-prefixed), matching the source data export where
# these were written as strings rather than numbers. A few cells are
# genuine numeric values (1, 0, 0.5, 0.71962616822429903) and are left
# as plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - plain text, no quote prefix needed (non-numeric strings
# are stored as shared strings automatically).
$ws.Range("A1").Value = "twoStep_rew"
$ws.Range("B1").Value = "twoStep_potent"
$ws.Range("C1").Value = "twoStep_percDead"

# Row 2
$ws.Range("A2").Value = "'0.8111111111111111"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "'0.3333333333333333"

# Row 3
$ws.Range("A3").Value = "'0.6304347826086957"
$ws.Range("B3").Value = "'0.7282608695652174"
$ws.Range("C3").Value = "'0.3333333333333333"

# Row 4
$ws.Range("A4").Value = "'0.42105263157894735"
$ws.Range("B4").Value = "'0.6491228070175439"
$ws.Range("C4").Value = 0.5

# Row 5
$ws.Range("A5").Value = "'0.5925925925925926"
$ws.Range("B5").Value = "'0.8024691358024691"
$ws.Range("C5").Value = 0.5

# Row 6
$ws.Range("A6").Value = "'0.7008547008547008"
$ws.Range("B6").Value = "'0.7777777777777778"
$ws.Range("C6").Value = 0.5

# Row 7
$ws.Range("A7").Value = "'0.8709677419354839"
$ws.Range("B7").Value = "'0.8709677419354839"
$ws.Range("C7").Value = 0

# Row 8
$ws.Range("A8").Value = "'0.7543859649122807"
$ws.Range("B8").Value = "'0.9035087719298246"
$ws.Range("C8").Value = "'0.3333333333333333"

# Row 9
$ws.Range("A9").Value = "'0.45555555555555555"
$ws.Range("B9").Value = "'0.8222222222222222"
$ws.Range("C9").Value = 0.5

# Row 10
$ws.Range("A10").Value = "'0.21153846153846154"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "'0.6666666666666666"

# Row 11
$ws.Range("A11").Value = "'0.8387096774193549"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "'0.3333333333333333"

# Row 12
$ws.Range("A12").Value = "'0.8245614035087719"
$ws.Range("B12").Value = "'0.9473684210526315"
$ws.Range("C12").Value = "'0.3333333333333333"

# Row 13
$ws.Range("A13").Value = "'0.6470588235294118"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0.5

# Row 14
$ws.Range("A14").Value = "'0.8021978021978022"
$ws.Range("B14").Value = "'0.9230769230769231"
$ws.Range("C14").Value = 0.5

# Row 15
$ws.Range("A15").Value = "'0.35514018691588783"
$ws.Range("B15").Value = 0.71962616822429903
$ws.Range("C15").Value = "'0.6666666666666666"

# Match the author's final selection (cell H6) recorded in the saved view.
$ws.Range("H6").Select()
